$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A; this shifts B:F left to A:E and drops the
# old column A along with its cell formatting (bold/bordered style).
$ws.Range("A1:A5").EntireColumn.Delete()
